$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Apply updated values
$ws.Range("D2").Value = "29.295.56"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.871.41"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "0.7124"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "241.79"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.3107"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "0.07718"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").Value = "24.76"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "0.08392"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "1.889.10"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "5.229"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "0.7132"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "91.16"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "29.288.13"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "0.000008203"
$ws.Range("E17").Value = "  +4.50%  "
$ws.Range("D18").Value = "5.941"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "243.61"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "2.125.17"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "7.881"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "0.1619"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "163.83"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "9.016"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "18.51"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "1.509"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "4.406"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "1.304"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").Value = "4.284"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("D33").Value = "0.05182"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7756"
$ws.Range("E34").Value = "  +6.69%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.918"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "2.711"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "1.156.86"
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("D41").Value = "6.399"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("D42").Value = "0.8922"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").Value = "73.25"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "103.63"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").Value = "2.020.03"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").Value = "1.800"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "0.5192"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "9.393"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "0.4303"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "7.052"
$ws.Range("E51").Value = "  -0.26%  "
